$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.521.92'
$ws.Range('D3').Value = '2.699.73'
$ws.Range('E3').Value = '  +2.20%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.48'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.66%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D9').Value = '2.698.24'
$ws.Range('E9').Value = '  +2.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.140'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.31%  '
$ws.Range('E11').Value = '  -0.29%  '
$ws.Range('E12').Value = '  +1.07%  '
$ws.Range('E13').Value = '  +2.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.31'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.41%  '
$ws.Range('D15').Value = '3.190.00'
$ws.Range('E15').Value = '  +2.15%  '
$ws.Range('E16').Value = '  -0.84%  '
$ws.Range('D17').Value = '68.506.40'
$ws.Range('E17').Value = '  +0.53%  '
$ws.Range('D18').Value = '2.695.91'
$ws.Range('E18').Value = '  +1.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.83'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.20%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.64'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.42%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '365.01'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.68%  '
$ws.Range('E22').Value = '  +2.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.89'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.31%  '
$ws.Range('E24').Value = '  +2.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '74.28'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.36%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.83'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.00%  '
$ws.Range('D28').Value = '2.832.45'
$ws.Range('E28').Value = '  +1.82%  '
$ws.Range('E29').Value = '  +1.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '586.83'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.73%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.01'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.22'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.31%  '
$ws.Range('E33').Value = '  +2.70%  '
$ws.Range('E34').Value = '  +5.17%  '
$ws.Range('E35').Value = '  +3.80%  '
$ws.Range('E36').Value = '  +6.05%  '
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '160.88'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('E39').Value = '  +1.02%  '
$ws.Range('E40').Value = '  +2.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.39'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.63%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.99'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.10%  '
$ws.Range('E44').Value = '  +3.28%  '
$ws.Range('D46').Value = '0.0₆0317'
$ws.Range('E46').Value = '  -4.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '157.41'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.68%  '
$ws.Range('E48').Value = '  +5.91%  '
$ws.Range('E49').Value = '  +5.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.604'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.05'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.64%  '
